$d = $word.ActiveDocument

# The target range spans from the start of the "Analog write" paragraph
# (currently paragraph #31) through the end of the empty bookmark-holding
# paragraph right after "Digital write" (currently paragraph #34).
$pStart = $d.Paragraphs.Item(31)
$pEnd = $d.Paragraphs.Item(34)
$target = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$newBody = @'
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Analog</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>write</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Digital </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>read</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Digital </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>write</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:jc w:val="both"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Tensión real: 250 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>mV</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Tension</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Medida: 0.16 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>u.a</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tensión generador: 1V</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Tensión real: 475 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>mV</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tensión generador: 2V</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Tensión medida: 0.3068 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>u.a</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tensión generador: 2.5V</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tensión medida: 0.3816</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Tensión </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>osci</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">: 600 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>mV</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Tension</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> generador: 3 V</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Tensión </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>emdida</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>: 0.4579</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Tensión </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>osci</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">: 700 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>mV</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tensión generador: 4V</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tensión medida: 0.610</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Tensión </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>osci</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>: 936 mV</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p/>
'@

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $newBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)

Write-Output ("Paragraph count after edit: " + $d.Paragraphs.Count)
